$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header E1 from 1560 to 1570 (PMD GUI parameter from user)
$ws.Range("E1").Value = 1570

# Append PMD GUI computed rows 17-26 (wavelengths 1561-1570 nm)
$ws.Range("A17").Value = "1561"
$ws.Range("B17").Value = "8.5642638734649905E-8"
$ws.Range("C17").Value = "2.07713909152692E-8"

$ws.Range("A18").Value = "1562"
$ws.Range("B18").Value = "8.5752401842750003E-8"
$ws.Range("C18").Value = "2.07980123792982E-8"

$ws.Range("A19").Value = "1563"
$ws.Range("B19").Value = "8.5862235244228898E-8"
$ws.Range("C19").Value = "2.0824650891975699E-8"

$ws.Range("A20").Value = "1564"
$ws.Range("B20").Value = "8.5972138939084301E-8"
$ws.Range("C20").Value = "2.0851306453301299E-8"

$ws.Range("A21").Value = "1565"
$ws.Range("B21").Value = "8.6082112927310401E-8"
$ws.Range("C21").Value = "2.0877979063273401E-8"

$ws.Range("A22").Value = "1566"
$ws.Range("B22").Value = "8.6192157208910905E-8"
$ws.Range("C22").Value = "2.0904668721893101E-8"

$ws.Range("A23").Value = "1567"
$ws.Range("B23").Value = "8.6302271783892601E-8"
$ws.Range("C23").Value = "2.09313754291619E-8"

$ws.Range("A24").Value = "1568"
$ws.Range("B24").Value = "8.6412456652258905E-8"
$ws.Range("C24").Value = "2.0958099185080699E-8"

$ws.Range("A25").Value = "1569"
$ws.Range("B25").Value = "8.6522711814003E-8"
$ws.Range("C25").Value = "2.09848399896477E-8"

$ws.Range("A26").Value = "1570"
$ws.Range("B26").Value = "8.6633037269121407E-8"
$ws.Range("C26").Value = "2.10115978428623E-8"
